$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sentiment-count table grew from 3 categories to 5 (CSV re-export produced
# more rows). Insert two formatted rows at 4 and 5 so the new "Very Positive"
# and "Negative" rows pick up the same (inherited) cell style as the existing
# rows, then drop the extra duplicate row introduced by the two inserts.
$ws.Rows("4:4").Insert(-4121)
$ws.Rows("5:5").Insert(-4121)
$ws.Rows("7:7").Delete()

# Column A labels (sentiment categories)
$ws.Range("A1").Value = "Very Negative"
$ws.Range("A2").Value = "Neutral"
$ws.Range("A3").Value = "Positive"
$ws.Range("A4").Value = "Very Positive"
$ws.Range("A5").Value = "Negative"

# Column B counts
$ws.Range("B1").Value = 6
$ws.Range("B2").Value = 18
$ws.Range("B3").Value = 9
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 61
